# Apply French translation updates to en_fr_maj.xlsx (Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "(Pour vous déconnecter, fermer l'app.)"
$ws.Range("A11").Value = "ACORN Participating Countries"
$ws.Range("B11").Value = "Pays Participants au Projet ACORN"
$ws.Range("A12").Value = "All 'orgname' are provided."
$ws.Range("B12").Value = "Tous les 'orgname' sont fournis."
$ws.Range("A13").Value = "All 'patid' are provided."
$ws.Range("B13").Value = "Tous les 'patid' sont fournis."
$ws.Range("A14").Value = "All 'specdate' are provided."
$ws.Range("B14").Value = "Tous les 'specdate' sont fournis."
$ws.Range("A15").Value = "All 'specdate' are today or before today."
$ws.Range("B15").Value = "Tous les 'specdate' sont égales ou antérieures à aujourd'hui."
$ws.Range("A16").Value = "All 'specgroup' are provided."
$ws.Range("B16").Value = "Tous les 'specgroup' sont fournis."
$ws.Range("A17").Value = "All 'specid' are provided."
$ws.Range("B17").Value = "Tous les 'specid' sont fournis."
$ws.Range("A18").Value = "All dates of enrolment for HAI patients have a matching date in the HAI survey dataset"
$ws.Range("B18").Value = "Toutes les dates d'enrôlements de patients avec une HAI se trouvent dans le jeu de données d'études de HAI."
$ws.Range("A19").Value = "All Other Organisms"
$ws.Range("B19").Value = "Tous les Autres Organismes"
$ws.Range("A20").Value = "All valid records have an ACORN ID."
$ws.Range("B20").Value = "Tous les enregistrements valides ont une ID ACORN."
$ws.Range("A21").Value = "AMR"
$ws.Range("B21").Value = "AMR"
$ws.Range("A22").Value = "and generate enrolment log."
$ws.Range("B22").Value = "et générer un fichier de suivi des recrutements."
$ws.Range("A23").Value = "Attempting to connect."
$ws.Range("B23").Value = "Tentative de connection."
$ws.Range("A24").Value = "Blood culture collected within 24 hours of admission (CAI) / symptom onset (HAI)"
$ws.Range("B24").Value = "Hémoculture collectée dans les 24 heures d'admission (CAI) / au début des symptômes (HAI)"
$ws.Range("A25").Value = "Blood Culture Contaminants"
$ws.Range("B25").Value = "Contaminants d'Hémoculture"
$ws.Range("A26").Value = "Bloodstream Infection (BSI)"
$ws.Range("B26").Value = "Infection du Sang (BSI)"
$ws.Range("A27").Value = "Calculated age is consistent with 'Age Category'"
$ws.Range("B27").Value = "L'âge calculé est cohérent avec la catégorie d'âge."
$ws.Range("A28").Value = "Calculated age isn't always consistent with 'Age Category'"
$ws.Range("B28").Value = "L'âge calculé n'est pas toujours cohérent avec la catégorie d'âge."
$ws.Range("A29").Value = "Cancel"
$ws.Range("B29").Value = "Annuler"
$ws.Range("A30").Value = "Care should be taken when interpreting rates and AMR profiles where there are small numbers of cases or bacterial isolates: point estimates may be unreliable."
$ws.Range("B30").Value = "L'interprétation des taux et profils AMR doit prêter à caution lorsqu'en présence d'un petit nombre de cas ou d'isolats : les estimations sont alors peu fiable."
$ws.Range("A31").Value = "Clinical and day-28 outcomes are consistent."
$ws.Range("B31").Value = "Les résultats cliniques et du jour 28 sont cohérents."
$ws.Range("A32").Value = "Clinical and day-28 outcomes aren't consistent for some dead patients."
$ws.Range("B32").Value = "Pour certains patients décédés, les résultats cliniques et du jour 28 ne sont pas cohérents."
$ws.Range("A33").Value = "Clinical Outcome"
$ws.Range("B33").Value = "Résultat Clinique"
$ws.Range("A34").Value = "Clinical Outcome Status:"
$ws.Range("B34").Value = "Statut des résultats cliniques:"
$ws.Range("A35").Value = "Co-resistances"
$ws.Range("B35").Value = "Co-résistances"
$ws.Range("A36").Value = "Combine Susceptible + Intermediate"
$ws.Range("B36").Value = "Combiner Susceptible + Intermédiaire"
$ws.Range("A37").Value = "Consider saving .acorn file on the cloud for additional security."
$ws.Range("B37").Value = "Pour plus de sécurité, envisagez une sauvegarde des données sur un serveur."
$ws.Range("A38").Value = "Contains names of organisms before and after mapping."
$ws.Range("B38").Value = "Contient les noms des organismes avant et après cartorgraphie."
$ws.Range("A39").Value = "Couldn't connect to server. Please check internet access."
$ws.Range("B39").Value = "Il n'est pas possible de se connecter au server. Veuillez vérifier la connection à internet"
$ws.Range("A40").Value = "Critical errors with clinical data."
$ws.Range("B40").Value = "Erreurs critiques avec les données cliniques."
$ws.Range("A41").Value = "Culture results per specimen type"
$ws.Range("B41").Value = "Résultats de culture par type de spécimen"
$ws.Range("A42").Value = "Data Management"
$ws.Range("B42").Value = "Gestion des Données"
$ws.Range("A43").Value = "Date of Enrolment"
$ws.Range("B43").Value = "Date d'Enrôlement"
$ws.Range("A44").Value = "Day 28"
$ws.Range("B44").Value = "Jour 28"
$ws.Range("A45").Value = "Day 28 Status:"
$ws.Range("B45").Value = "Statut au jour 28:"
$ws.Range("A46").Value = "Diagnosis at Enrolment"
$ws.Range("B46").Value = "Diagnostics à l'Enrôlement"
$ws.Range("A47").Value = "Dismiss"
$ws.Range("B47").Value = "Annuler"
$ws.Range("A48").Value = "Distribution of Enrolments"
$ws.Range("B48").Value = "Distribution des Enrôlements"
$ws.Range("A49").Value = "Download Enrolment Log (.xlsx)"
$ws.Range("B49").Value = "Télécharger un fichier de suivi des recrutements (.xlsx)"
$ws.Range("A50").Value = "Download Lab Log (.xlsx)"
$ws.Range("B50").Value = "Télécharger un log du traitement des donnés de lab (.xlsx)"
$ws.Range("B65").Value = "Générer un fichier .acorn depuis les données cliniques et de lab"
$ws.Range("B68").Value = "Obtenir la dernière version stable"
$ws.Range("A70").Value = "HAI point prevalence by "
$ws.Range("B70").Value = "TBT"
$ws.Range("B72").Value = "Les barres horizontales indiquent la taille d'un ensemble de résultats SR tandis que les barres verticales indiquent le nombre d'isolats résistants pour l'antibiotique correspondant."
$ws.Range("B73").Value = "Information sur le fichier .acorn chargé."
$ws.Range("B76").Value = "Problème détecté avec les données REDCap. Merci de contacter l'équipe ACORN. Jusqu'à résolution, seuls les fichiers .acorn existants peuvent être utilisés."
$ws.Range("B84").Value = "Charger le fichier .acorn depuis le nuage"
$ws.Range("B85").Value = "Charger le fichier .acorn localement"
$ws.Range("B100").Value = "Seuls les isolats qui ont été testés contre tous les médicaments sont inclus dans le graphique."
$ws.Range("A110").Value = "Remove 'Not Cultured' specimens"
$ws.Range("B110").Value = "Supprimer les spécimens 'Not Cultured'"
$ws.Range("A111").Value = "Remove blood culture contaminants from the following visualizations"
$ws.Range("B111").Value = "Élimine les contaminants d'hémoculture des visualisations ci-dessous"
$ws.Range("A112").Value = "Reset Enrolments Filters"
$ws.Range("B112").Value = "Réinitialiser les Filtres sur les Enrôlements"
$ws.Range("A113").Value = "Resistance to 3rd gen. Cephalosporins Over Time"
$ws.Range("B113").Value = "Evolution de la Résistance aux 3rd gen. Cephalosporins"
$ws.Range("A114").Value = "Resistance to Carbapenems Over Time"
$ws.Range("B114").Value = "Evolution de la Résistance aux Carbapenems"
$ws.Range("A115").Value = "Resistance to Fluoroquinolones Over Time"
$ws.Range("B115").Value = "Evolution de la Résistance aux Fluoroquinolones"
$ws.Range("A116").Value = "Resistance to Oxacillin Over Time"
$ws.Range("B116").Value = "Evolution de la Résistance aux Oxacillin"
$ws.Range("A117").Value = "Resistance to Penicillin G - meningitis Over Time"
$ws.Range("B117").Value = "Evolution de la Résistance aux Penicillin G - meningitis"
$ws.Range("A118").Value = "Resistance to Penicillin G Over Time"
$ws.Range("B118").Value = "Evolution de la Résistance aux Penicillin G"
$ws.Range("A119").Value = "Retriving data from REDCap server."
$ws.Range("B119").Value = "Récupération des données depuis le serveur REDcap."
$ws.Range("A120").Value = "Save .acorn file"
$ws.Range("B120").Value = "Sauvegarder un ficher .acorn"
$ws.Range("A121").Value = "Save acorn data"
$ws.Range("B121").Value = "Sauvegarder des données acorn"
$ws.Range("A122").Value = "Save on Server"
$ws.Range("B122").Value = "Sauvegarder sur serveur"
$ws.Range("A123").Value = "See Breakdown by Ward"
$ws.Range("B123").Value = "Montrer la Répartition par Service"
$ws.Range("A124").Value = "See by Week"
$ws.Range("B124").Value = "Montrer par Semaine"
$ws.Range("B125").Value = "Montrer les combinaisons d'antibiotiques."
$ws.Range("B127").Value = "Evaluation SIR"
$ws.Range("B133").Value = "Certains enregistrements ont un identifiant ACORN manquant. Ces enregistrements ont été supprimés."
$ws.Range("B135").Value = "Spécimens"
$ws.Range("B144").Value = "Susceptible & Intermédiaire sont toujours combinés dans cette visualisation des co-résistances."
$ws.Range("B146").Value = "Les « identifiants de patient » suivants sont des cas atypiques (un HCAI/CAI avec HAI précoce mais sans chevauchement) :"
$ws.Range("B148").Value = "Le jeu de données REDCap est vide/au mauvais format. Veuillez contacter l'assistance ACORN."
$ws.Range("B149").Value = "Le jeu de données REDCap est au bon format."
$ws.Range("B150").Value = "Il y a des suivis à J28 effectués avant la date prévue à J28."
$ws.Range("B151").Value = "Il existe plusieurs F02 avec un ID ACORN, une date d'admission et une date d'enrôlement identiques."
$ws.Range("B152").Value = "Il n'y a pas de cas atypique (un HCAI/CAI avec HAI précoce mais pas de chevauchement)."
$ws.Range("B153").Value = "Il n'y a pas de suivi à J28 effectué avant la date prévue à J28."
$ws.Range("B155").Value = "Il n'existe pas de F02 avec un ID ACORN, une date d'admission et une date d'enrôlement identiques."
$ws.Range("B167").Value = "Updated Charlson Comorbidity Index (uCCI)"
$ws.Range("B172").Value = "Le dictionnaire des données de lab ne peut pas être téléchargé. Merci de contacter l'équipe ACORN."
$ws.Range("B177").Value = "Vous utilisez le tableau de bord ACORN"
$ws.Range("B178").Value = "Vous pouvez vérifier ici s'il s'agit de la dernière version de production."
$ws.Range("B179").Value = "Votre tableau de bord ACORN est à jour"
